$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove all sheets except the first one (PG), which holds the real data.
$namesToRemove = @("SG", "PF", "SF", "C", "All", "STATS", "Sheet8")
foreach ($name in $namesToRemove) {
    $wb.Worksheets.Item($name).Delete()
}

# Rename the remaining sheet to "ALL PLAYERS".
$wb.Worksheets.Item("PG").Name = "ALL PLAYERS"

$excel.DisplayAlerts = $true
